# Generate Report for Handback
# - Update the "Status" column (now that localized content is back in sync
#   with en-US, the old "Ready for handoff" label becomes "Handed back: in
#   sync with en-US") on the Overview sheet and on each language sheet.
# - Populate the "Latest Target File" (F) and "Latest Handback File" (G)
#   columns -- together with their hyperlinks -- on the zh-cn and de-de
#   sheets, now that the handback has actually happened.
# - Record the real handback timestamp in the "Latest Handback DateTime"
#   (H) column, replacing the zero-value placeholder.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---- Status column updates (shared string used across all 3 sheets) ----
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---- zh-cn sheet: add Latest Target File / Latest Handback File columns ----
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a27a8fe05c0f05dcb120124cbd54e19d004bd71f/e2e/28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.md",
    "",
    "",
    "28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d54fed170564b5d6cc3a1a0d3e3ce1bfa8169156/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.48ff7a36ebd8f9d35dac1bf4043e4df4631bd037.zh-cn.xlf",
    "",
    "",
    "28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.48ff7a36ebd8f9d35dac1bf4043e4df4631bd037.zh-cn.xlf"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a27a8fe05c0f05dcb120124cbd54e19d004bd71f/e2e/4503c4c5-86d4-4d3b-9b35-781577df6db2.md",
    "",
    "",
    "4503c4c5-86d4-4d3b-9b35-781577df6db2.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d54fed170564b5d6cc3a1a0d3e3ce1bfa8169156/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/4503c4c5-86d4-4d3b-9b35-781577df6db2.12da529969af735f7627f8fc7f7441ae729caeb6.zh-cn.xlf",
    "",
    "",
    "4503c4c5-86d4-4d3b-9b35-781577df6db2.12da529969af735f7627f8fc7f7441ae729caeb6.zh-cn.xlf"
)

# Latest Handback DateTime for zh-cn
$wsZhCn.Range("H2").Value = "2016-03-14 09:50:55"
$wsZhCn.Range("H3").Value = "2016-03-14 09:50:55"

# ---- de-de sheet: add Latest Target File / Latest Handback File columns ----
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a27a8fe05c0f05dcb120124cbd54e19d004bd71f/e2e/28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.md",
    "",
    "",
    "28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b018b80294ea185c1088ba8f7fe3f1732c8ef633/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.48ff7a36ebd8f9d35dac1bf4043e4df4631bd037.de-de.xlf",
    "",
    "",
    "28e66ae9-a8ed-41a3-9b13-d6847b95ad0f.48ff7a36ebd8f9d35dac1bf4043e4df4631bd037.de-de.xlf"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a27a8fe05c0f05dcb120124cbd54e19d004bd71f/e2e/4503c4c5-86d4-4d3b-9b35-781577df6db2.md",
    "",
    "",
    "4503c4c5-86d4-4d3b-9b35-781577df6db2.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b018b80294ea185c1088ba8f7fe3f1732c8ef633/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/4503c4c5-86d4-4d3b-9b35-781577df6db2.12da529969af735f7627f8fc7f7441ae729caeb6.de-de.xlf",
    "",
    "",
    "4503c4c5-86d4-4d3b-9b35-781577df6db2.12da529969af735f7627f8fc7f7441ae729caeb6.de-de.xlf"
)

# Latest Handback DateTime for de-de
$wsDeDe.Range("H2").Value = "2016-03-14 09:51:09"
$wsDeDe.Range("H3").Value = "2016-03-14 09:51:09"
